{"js": "// Replace the two-digit multiplication problems in the document's table\n// with a new set of problems, preserving all formatting. Each original\n// expression (e.g. \"97\u00d757=\") is unique in the document, so we can find\n// each one and replace it in place.\nconst replacements = [\n  [\"97\u00d757=\", \"48\u00d783=\"],\n  [\"27\u00d727=\", \"93\u00d778=\"],\n  [\"75\u00d790=\", \"51\u00d776=\"],\n  [\"48\u00d755=\", \"64\u00d725=\"],\n  [\"13\u00d789=\", \"33\u00d745=\"],\n  [\"95\u00d762=\", \"81\u00d732=\"],\n  [\"71\u00d746=\", \"97\u00d711=\"],\n  [\"68\u00d731=\", \"97\u00d778=\"],\n  [\"94\u00d799=\", \"89\u00d785=\"],\n  [\"15\u00d787=\", \"13\u00d790=\"],\n  [\"96\u00d766=\", \"29\u00d741=\"],\n  [\"27\u00d735=\", \"66\u00d747=\"],\n  [\"11\u00d786=\", \"14\u00d788=\"],\n  [\"48\u00d761=\", \"21\u00d793=\"],\n  [\"69\u00d793=\", \"37\u00d711=\"],\n  [\"42\u00d764=\", \"97\u00d754=\"],\n  [\"81\u00d782=\", \"86\u00d760=\"],\n  [\"39\u00d744=\", \"77\u00d779=\"],\n  [\"98\u00d763=\", \"71\u00d790=\"],\n  [\"19\u00d717=\", \"37\u00d793=\"],\n  [\"75\u00d717=\", \"82\u00d774=\"],\n  [\"88\u00d722=\", \"42\u00d754=\"],\n  [\"23\u00d719=\", \"76\u00d777=\"],\n  [\"33\u00d773=\", \"99\u00d752=\"],\n  [\"81\u00d765=\", \"74\u00d762=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the document's table\n# with a new set of problems, preserving all formatting. Each original\n# expression (e.g. \"97x57=\") is unique in the document, so Find/Replace\n# on the whole story finds exactly one match per pair.\n\n$d = $word.ActiveDocument\n$mult = [char]215  # \"x\" multiplication sign (U+00D7)\n\n$replacements = @(\n    @(\"97${mult}57=\", \"48${mult}83=\"),\n    @(\"27${mult}27=\", \"93${mult}78=\"),\n    @(\"75${mult}90=\", \"51${mult}76=\"),\n    @(\"48${mult}55=\", \"64${mult}25=\"),\n    @(\"13${mult}89=\", \"33${mult}45=\"),\n    @(\"95${mult}62=\", \"81${mult}32=\"),\n    @(\"71${mult}46=\", \"97${mult}11=\"),\n    @(\"68${mult}31=\", \"97${mult}78=\"),\n    @(\"94${mult}99=\", \"89${mult}85=\"),\n    @(\"15${mult}87=\", \"13${mult}90=\"),\n    @(\"96${mult}66=\", \"29${mult}41=\"),\n    @(\"27${mult}35=\", \"66${mult}47=\"),\n    @(\"11${mult}86=\", \"14${mult}88=\"),\n    @(\"48${mult}61=\", \"21${mult}93=\"),\n    @(\"69${mult}93=\", \"37${mult}11=\"),\n    @(\"42${mult}64=\", \"97${mult}54=\"),\n    @(\"81${mult}82=\", \"86${mult}60=\"),\n    @(\"39${mult}44=\", \"77${mult}79=\"),\n    @(\"98${mult}63=\", \"71${mult}90=\"),\n    @(\"19${mult}17=\", \"37${mult}93=\"),\n    @(\"75${mult}17=\", \"82${mult}74=\"),\n    @(\"88${mult}22=\", \"42${mult}54=\"),\n    @(\"23${mult}19=\", \"76${mult}77=\"),\n    @(\"33${mult}73=\", \"99${mult}52=\"),\n    @(\"81${mult}65=\", \"74${mult}62=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
